$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2:D51').NumberFormat = '@'
$ws.Range('E2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '56.050.31'
$ws.Range('E2').Value = '  +5.45%  '
$ws.Range('D3').Value = '2.517.26'
$ws.Range('E3').Value = '  +6.24%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '491.43'
$ws.Range('E5').Value = '  +8.02%  '
$ws.Range('D6').Value = '144.39'
$ws.Range('E6').Value = '  +12.82%  '
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '0.515'
$ws.Range('E8').Value = '  +7.56%  '
$ws.Range('D9').Value = '2.536.48'
$ws.Range('E9').Value = '  +6.40%  '
$ws.Range('D10').Value = '0.0987'
$ws.Range('E10').Value = '  +6.36%  '
$ws.Range('D11').Value = '5.62'
$ws.Range('E11').Value = '  +6.16%  '
$ws.Range('D12').Value = '0.334'
$ws.Range('E12').Value = '  +6.38%  '
$ws.Range('E13').Value = '  +1.47%  '
$ws.Range('D14').Value = '2.947.88'
$ws.Range('E14').Value = '  +6.24%  '
$ws.Range('D15').Value = '56.041.04'
$ws.Range('E15').Value = '  +5.38%  '
$ws.Range('D16').Value = '21.11'
$ws.Range('E16').Value = '  +8.89%  '
$ws.Range('D17').Value = '0.0000137'
$ws.Range('E17').Value = '  +11.23%  '
$ws.Range('D18').Value = '2.512.88'
$ws.Range('E18').Value = '  +6.53%  '
$ws.Range('D19').Value = '4.45'
$ws.Range('E19').Value = '  +7.44%  '
$ws.Range('D20').Value = '10.25'
$ws.Range('E20').Value = '  +12.30%  '
$ws.Range('D21').Value = '322.31'
$ws.Range('E21').Value = '  +5.56%  '
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -2.08%  '
$ws.Range('D23').Value = '5.86'
$ws.Range('E23').Value = '  +9.99%  '
$ws.Range('D24').Value = '58.33'
$ws.Range('E24').Value = '  +6.05%  '
$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D25').Value = '0.168'
$ws.Range('E25').Value = '  +13.86%  '
$ws.Range('D26').Value = '0.412'
$ws.Range('E26').Value = '  +8.81%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  +0.23%  '
$ws.Range('D28').Value = '2.581.24'
$ws.Range('E28').Value = '  +6.67%  '
$ws.Range('D29').Value = '7.50'
$ws.Range('E29').Value = '  +6.59%  '
$ws.Range('D30').Value = '0.0₃0793'
$ws.Range('E30').Value = '  +15.13%  '
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  +0.39%  '
$ws.Range('D32').Value = '149.27'
$ws.Range('E32').Value = '  +1.48%  '
$ws.Range('D33').Value = '18.51'
$ws.Range('E33').Value = '  +6.35%  '
$ws.Range('D34').Value = '1.51'
$ws.Range('E34').Value = '  +9.69%  '
$ws.Range('D35').Value = '5.28'
$ws.Range('E35').Value = '  +7.39%  '
$ws.Range('D36').Value = '1.16'
$ws.Range('E36').Value = '  +13.26%  '
$ws.Range('D37').Value = '3.73'
$ws.Range('E37').Value = '  +8.14%  '
$ws.Range('D38').Value = '0.871'
$ws.Range('E38').Value = '  +5.70%  '
$ws.Range('D39').Value = '34.17'
$ws.Range('E39').Value = '  +3.86%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').Value = '3.52'
$ws.Range('E40').Value = '  +8.35%  '
$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D41').Value = '0.621'
$ws.Range('E41').Value = '  +5.30%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').Value = '0.0556'
$ws.Range('E42').Value = '  +7.21%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = '0.993'
$ws.Range('E43').Value = '  -0.22%  '
$ws.Range('D44').Value = '1.33'
$ws.Range('E44').Value = '  +11.62%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = '4.88'
$ws.Range('E45').Value = '  +16.40%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').Value = '268.23'
$ws.Range('E46').Value = '  +32.23%  '
$ws.Range('B47').Value = 'WhiteBITCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D47').Value = '10.17'
$ws.Range('E47').Value = '  +0.60%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '0.0227'
$ws.Range('E48').Value = '  +6.94%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = '0.0905'
$ws.Range('E49').Value = '  +6.60%  '
$ws.Range('D50').Value = '1.962.66'
$ws.Range('E50').Value = '  +1.17%  '
$ws.Range('D51').Value = '17.77'
$ws.Range('E51').Value = '  +8.38%  '
